$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append row 4 with the new trip/facility record. All cells in this sheet
# are stored as text, including values that look numeric (e.g. "12"), and
# A4 is an explicit empty-string cell (not a blank cell).

# A4: explicit empty string. A leading apostrophe forces text entry with an
# empty value; resetting the style afterwards drops the quote-prefix flag
# so the cell ends up with the default (unstyled) text cell the diff wants.
$ws.Range("A4").Value = "'"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "أحمد شريم"

# C4: "12" looks numeric, so force text via NumberFormat, then reset the
# style so no extra number-format style sticks to the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "12"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "الجزائري"
$ws.Range("E4").Value = "الرحلة 1"
$ws.Range("F4").Value = "C1"
$ws.Range("G4").Value = "UNICEF"
$ws.Range("H4").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٤١:٤٩ م"
